$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "264.95"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.74"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.242"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06161"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.592"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.698"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.357"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1608"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08180"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03392"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09252"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.912"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001713"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04790"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006292"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001100"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3341"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1252"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04639"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006971"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1137"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003131"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006164"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7780"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2041"
